# Updated cryptos list - refresh Price/Volume(1h) values and shift rows for new USDe entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing formatting like trailing zeros)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '63.519.95'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '3.401.08'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '567.25'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').Value = '156.47'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.403.27'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('D9').Value = '0.568'
$ws.Range('E9').Value = '  -8.37%  '
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('E11').Value = '  -3.60%  '
$ws.Range('D12').Value = '0.423'
$ws.Range('E12').Value = '  -4.34%  '
$ws.Range('D13').Value = '3.987.57'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '26.83'
$ws.Range('E15').Value = '  -4.02%  '
$ws.Range('E16').Value = '  -9.23%  '
$ws.Range('D17').Value = '63.585.07'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').Value = '3.392.34'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('E19').Value = '  -4.72%  '
$ws.Range('D20').Value = '13.49'
$ws.Range('E20').Value = '  -3.46%  '
$ws.Range('D21').Value = '383.03'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').Value = '7.73'
$ws.Range('E22').Value = '  -3.48%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = '71.05'
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('E25').Value = '  -7.47%  '
$ws.Range('E26').Value = '  -3.98%  '
$ws.Range('D27').Value = '9.67'
$ws.Range('E27').Value = '  -6.09%  '
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  -2.98%  '
$ws.Range('E31').Value = '  -7.15%  '
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '22.80'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '6.90'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  -7.18%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '160.36'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = '0.840'
$ws.Range('E38').Value = '  +9.13%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.81'
$ws.Range('E39').Value = '  -4.50%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.813.21'
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '25.86'
$ws.Range('E41').Value = '  -3.08%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '42.92'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '0.0713'
$ws.Range('E43').Value = '  -6.55%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '6.39'
$ws.Range('E44').Value = '  -8.20%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '25.56'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '4.34'
$ws.Range('E46').Value = '  -6.16%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0302'
$ws.Range('E47').Value = '  -4.18%  '
$ws.Range('E48').Value = '  +7.82%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '324.45'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '1.02'
$ws.Range('E50').Value = '  -5.13%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.103'
$ws.Range('E51').Value = '  -5.69%  '
